$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F/G column values for rows 518-705 (AgTests / AgPosit corrections)
$ws.Range("F518").Value = 7379

$ws.Range("F533").Value = 11986

$ws.Range("F582").Value = 26364

$ws.Range("F597").Value = 29876

$ws.Range("F608").Value = 46654

$ws.Range("F623").Value = 15356

$ws.Range("F640").Value = 19893

$ws.Range("F646").Value = 36044

$ws.Range("F657").Value = 34076

$ws.Range("F664").Value = 26429

$ws.Range("F665").Value = 28240

$ws.Range("F670").Value = 52642

$ws.Range("F678").Value = 33842

$ws.Range("F679").Value = 29464

$ws.Range("F684").Value = 57098

$ws.Range("F688").Value = 32027

$ws.Range("F691").Value = 62239

$ws.Range("F692").Value = 41525

$ws.Range("F694").Value = 37350
$ws.Range("G694").Value = 2772

$ws.Range("F695").Value = 36631
$ws.Range("G695").Value = 3081

$ws.Range("F696").Value = 17575
$ws.Range("G696").Value = 2184

$ws.Range("F697").Value = 28304
$ws.Range("G697").Value = 2965

$ws.Range("F698").Value = 67909
$ws.Range("G698").Value = 5694

$ws.Range("F699").Value = 42530

$ws.Range("F700").Value = 42226
$ws.Range("G700").Value = 4122

$ws.Range("F701").Value = 40983
$ws.Range("G701").Value = 3762

$ws.Range("F702").Value = 35145
$ws.Range("G702").Value = 3795

$ws.Range("F703").Value = 16516
$ws.Range("G703").Value = 2519

$ws.Range("F704").Value = 23431
$ws.Range("G704").Value = 3412

$ws.Range("F705").Value = 50593
$ws.Range("G705").Value = 5484

# Complete row 706 (add AgTests / AgPosit values)
$ws.Range("F706").Value = 36455
$ws.Range("G706").Value = 4290

# Add new row 707
$ws.Range("A707").Value = 44601
$ws.Range("B707").Value = 1187002
$ws.Range("C707").Value = 35106
$ws.Range("D707").Value = 20160
$ws.Range("E707").Value = 18040
$ws.Range("F707").Value = 25017
$ws.Range("G707").Value = 2982
